# Update the LR-pair table (Fgf18-Fgfr2) with the new TPM-derived values.
# The sending/target cluster set changed: "ECs" is now a sending cluster
# (previously only a target cluster) and "Resolving-Mac" is no longer a
# target cluster, so every data row (A2:T13) is rewritten in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = New-Object 'object[,]' 12,20

# Row 2
$data[0,0] = "ECs"
$data[0,1] = "Fgf18"
$data[0,2] = "Fgfr2"
$data[0,3] = "ECs"
$data[0,4] = 1
$data[0,5] = 0.3333333333333333
$data[0,6] = 0.1601763333333333
$data[0,7] = 0.480529
$data[0,8] = 0.01412814675921196
$data[0,9] = 0.01412814675921196
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 1.414593
$data[0,13] = 4.243779
$data[0,14] = 0.3478871232761722
$data[0,15] = 0.3478871232761722
$data[0,16] = 0.226584319899
$data[0,17] = 2.039258879091
$data[0,18] = 0.004915000333285823
$data[0,19] = 0.004915000333285824

# Row 3
$data[1,0] = "ECs"
$data[1,1] = "Fgf18"
$data[1,2] = "Fgfr2"
$data[1,3] = "FAPs"
$data[1,4] = 1
$data[1,5] = 0.3333333333333333
$data[1,6] = 0.1601763333333333
$data[1,7] = 0.480529
$data[1,8] = 0.01412814675921196
$data[1,9] = 0.01412814675921196
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 2.337487333333333
$data[1,13] = 7.012461999999999
$data[1,14] = 0.5748520910875596
$data[1,15] = 0.5748520910875596
$data[1,16] = 0.3744101502664444
$data[1,17] = 3.369691352398
$data[1,18] = 0.008121594707724923
$data[1,19] = 0.008121594707724923

# Row 4
$data[2,0] = "ECs"
$data[2,1] = "Fgf18"
$data[2,2] = "Fgfr2"
$data[2,3] = "MuSCs"
$data[2,4] = 1
$data[2,5] = 0.3333333333333333
$data[2,6] = 0.1601763333333333
$data[2,7] = 0.480529
$data[2,8] = 0.01412814675921196
$data[2,9] = 0.01412814675921196
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 0.314161
$data[2,13] = 0.942483
$data[2,14] = 0.07726078563626818
$data[2,15] = 0.07726078563626819
$data[2,16] = 0.05032115705633333
$data[2,17] = 0.452890413507
$data[2,18] = 0.001091551718201212
$data[2,19] = 0.001091551718201212

# Row 5
$data[3,0] = "FAPs"
$data[3,1] = "Fgf18"
$data[3,2] = "Fgfr2"
$data[3,3] = "ECs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 10.019353
$data[3,7] = 30.058059
$data[3,8] = 0.8837441004581448
$data[3,9] = 0.8837441004581448
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 1.414593
$data[3,13] = 4.243779
$data[3,14] = 0.3478871232761722
$data[3,15] = 0.3478871232761722
$data[3,16] = 14.173306618329
$data[3,17] = 127.559759564961
$data[3,18] = 0.3074431928206725
$data[3,19] = 0.3074431928206726

# Row 6
$data[4,0] = "FAPs"
$data[4,1] = "Fgf18"
$data[4,2] = "Fgfr2"
$data[4,3] = "FAPs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 10.019353
$data[4,7] = 30.058059
$data[4,8] = 0.8837441004581448
$data[4,9] = 0.8837441004581448
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 2.337487333333333
$data[4,13] = 7.012461999999999
$data[4,14] = 0.5748520910875596
$data[4,15] = 0.5748520910875596
$data[4,16] = 23.42011072569533
$data[4,17] = 210.780996531258
$data[4,18] = 0.5080221441346588
$data[4,19] = 0.5080221441346588

# Row 7
$data[5,0] = "FAPs"
$data[5,1] = "Fgf18"
$data[5,2] = "Fgfr2"
$data[5,3] = "MuSCs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 10.019353
$data[5,7] = 30.058059
$data[5,8] = 0.8837441004581448
$data[5,9] = 0.8837441004581448
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 0.314161
$data[5,13] = 0.942483
$data[5,14] = 0.07726078563626818
$data[5,15] = 0.07726078563626819
$data[5,16] = 3.147689957833
$data[5,17] = 28.329209620497
$data[5,18] = 0.06827876350281337
$data[5,19] = 0.06827876350281338

# Row 8
$data[6,0] = "MuSCs"
$data[6,1] = "Fgf18"
$data[6,2] = "Fgfr2"
$data[6,3] = "ECs"
$data[6,4] = 2
$data[6,5] = 0.6666666666666666
$data[6,6] = 1.010195666666666
$data[6,7] = 3.030587
$data[6,8] = 0.08910300502687639
$data[6,9] = 0.0891030050268764
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 1.414593
$data[6,13] = 4.243779
$data[6,14] = 0.3478871232761722
$data[6,15] = 0.3478871232761722
$data[6,16] = 1.429015718697
$data[6,17] = 12.861141468273
$data[6,18] = 0.03099778809406233
$data[6,19] = 0.03099778809406234

# Row 9
$data[7,0] = "MuSCs"
$data[7,1] = "Fgf18"
$data[7,2] = "Fgfr2"
$data[7,3] = "FAPs"
$data[7,4] = 2
$data[7,5] = 0.6666666666666666
$data[7,6] = 1.010195666666666
$data[7,7] = 3.030587
$data[7,8] = 0.08910300502687639
$data[7,9] = 0.0891030050268764
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 2.337487333333333
$data[7,13] = 7.012461999999999
$data[7,14] = 0.5748520910875596
$data[7,15] = 0.5748520910875596
$data[7,16] = 2.361319575021555
$data[7,17] = 21.251876175194
$data[7,18] = 0.05122104876188523
$data[7,19] = 0.05122104876188524

# Row 10
$data[8,0] = "MuSCs"
$data[8,1] = "Fgf18"
$data[8,2] = "Fgfr2"
$data[8,3] = "MuSCs"
$data[8,4] = 2
$data[8,5] = 0.6666666666666666
$data[8,6] = 1.010195666666666
$data[8,7] = 3.030587
$data[8,8] = 0.08910300502687639
$data[8,9] = 0.0891030050268764
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 0.314161
$data[8,13] = 0.942483
$data[8,14] = 0.07726078563626818
$data[8,15] = 0.07726078563626819
$data[8,16] = 0.3173640808356666
$data[8,17] = 2.856276727521
$data[8,18] = 0.006884168170928823
$data[8,19] = 0.006884168170928824

# Row 11
$data[9,0] = "Resolving-Mac"
$data[9,1] = "Fgf18"
$data[9,2] = "Fgfr2"
$data[9,3] = "ECs"
$data[9,4] = 1
$data[9,5] = 0.3333333333333333
$data[9,6] = 0.1476666666666667
$data[9,7] = 0.443
$data[9,8] = 0.01302474775576687
$data[9,9] = 0.01302474775576687
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 1.414593
$data[9,13] = 4.243779
$data[9,14] = 0.3478871232761722
$data[9,15] = 0.3478871232761722
$data[9,16] = 0.208888233
$data[9,17] = 1.879994097
$data[9,18] = 0.004531142028151515
$data[9,19] = 0.004531142028151516

# Row 12
$data[10,0] = "Resolving-Mac"
$data[10,1] = "Fgf18"
$data[10,2] = "Fgfr2"
$data[10,3] = "FAPs"
$data[10,4] = 1
$data[10,5] = 0.3333333333333333
$data[10,6] = 0.1476666666666667
$data[10,7] = 0.443
$data[10,8] = 0.01302474775576687
$data[10,9] = 0.01302474775576687
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 2.337487333333333
$data[10,13] = 7.012461999999999
$data[10,14] = 0.5748520910875596
$data[10,15] = 0.5748520910875596
$data[10,16] = 0.3451689628888889
$data[10,17] = 3.106520666
$data[10,18] = 0.007487303483290584
$data[10,19] = 0.007487303483290584

# Row 13
$data[11,0] = "Resolving-Mac"
$data[11,1] = "Fgf18"
$data[11,2] = "Fgfr2"
$data[11,3] = "MuSCs"
$data[11,4] = 1
$data[11,5] = 0.3333333333333333
$data[11,6] = 0.1476666666666667
$data[11,7] = 0.443
$data[11,8] = 0.01302474775576687
$data[11,9] = 0.01302474775576687
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 0.314161
$data[11,13] = 0.942483
$data[11,14] = 0.07726078563626818
$data[11,15] = 0.07726078563626819
$data[11,16] = 0.04639110766666666
$data[11,17] = 0.417519969
$data[11,18] = 0.001006302244324769
$data[11,19] = 0.001006302244324769

$ws.Range("A2:T13").Value = $data

